$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1 / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 990
$ws1.Range("F5").Value = 7059
$ws1.Range("F7").Value = 915
$ws1.Range("F10").Value = 538
$ws1.Range("F16").Value = 2823
$ws1.Range("F18").Value = 23
$ws1.Range("F25").Value = 98
$ws1.Range("F26").Value = 146
$ws1.Range("F32").Value = 240
$ws1.Range("F33").Value = 359

# Sheet "全部类型" (sheetId 4 / sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 990
$ws4.Range("F9").Value = 7059
$ws4.Range("F11").Value = 915
$ws4.Range("F14").Value = 538
$ws4.Range("F21").Value = 2823
$ws4.Range("F23").Value = 23
$ws4.Range("F32").Value = 98
$ws4.Range("F33").Value = 146
$ws4.Range("F39").Value = 240
$ws4.Range("F40").Value = 359
